$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '63.188.49'
$ws.Range('E2').Value = '  -5.23%  '
Set-TextValue 'D3' '3.311.73'
$ws.Range('E3').Value = '  -5.89%  '
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '549.51'
$ws.Range('E5').Value = '  -2.27%  '
Set-TextValue 'D6' '170.20'
$ws.Range('E6').Value = '  -8.50%  '
Set-TextValue 'D7' '0.607'
$ws.Range('E7').Value = '  -5.07%  '
$ws.Range('E8').Value = '  -0.03%  '
Set-TextValue 'D9' '3.295.40'
$ws.Range('E9').Value = '  -6.26%  '
Set-TextValue 'D10' '0.614'
$ws.Range('E10').Value = '  -4.61%  '
Set-TextValue 'D11' '0.150'
$ws.Range('E11').Value = '  -4.38%  '
Set-TextValue 'D12' '52.82'
$ws.Range('E12').Value = '  -5.07%  '
Set-TextValue 'D13' '0.0000263'
$ws.Range('E13').Value = '  -5.93%  '
Set-TextValue 'D14' '8.87'
$ws.Range('E14').Value = '  -6.09%  '
Set-TextValue 'D15' '3.836.54'
$ws.Range('E15').Value = '  -5.75%  '
Set-TextValue 'D16' '0.117'
$ws.Range('E16').Value = '  -3.96%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D17' '17.69'
$ws.Range('E17').Value = '  -5.44%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '3.291.28'
$ws.Range('E18').Value = '  -6.26%  '
Set-TextValue 'D19' '11.59'
$ws.Range('E19').Value = '  -4.55%  '
Set-TextValue 'D20' '63.102.41'
$ws.Range('E20').Value = '  -5.29%  '
Set-TextValue 'D21' '0.963'
$ws.Range('E21').Value = '  -4.26%  '
Set-TextValue 'D22' '403.11'
$ws.Range('E22').Value = '  -4.09%  '
Set-TextValue 'D23' '4.04'
$ws.Range('E23').Value = '  -1.54%  '
Set-TextValue 'D24' '4.24'
$ws.Range('E24').Value = '  +1.79%  '
Set-TextValue 'D25' '82.41'
$ws.Range('E25').Value = '  -4.44%  '
Set-TextValue 'D26' '13.12'
$ws.Range('E26').Value = '  +5.48%  '
Set-TextValue 'D27' '10.53'
$ws.Range('E27').Value = '  -4.29%  '
Set-TextValue 'D28' '2.71'
$ws.Range('E28').Value = '  -7.55%  '
Set-TextValue 'D29' '8.54'
$ws.Range('E29').Value = '  -7.69%  '
Set-TextValue 'D30' '28.90'
$ws.Range('E30').Value = '  -5.25%  '
Set-TextValue 'D31' '6.37'
$ws.Range('E31').Value = '  -5.53%  '
Set-TextValue 'D32' '11.25'
$ws.Range('E32').Value = '  -5.86%  '
Set-TextValue 'D33' '568.63'
$ws.Range('E33').Value = '  -8.15%  '
$ws.Range('E34').Value = '  -5.98%  '
Set-TextValue 'D35' '57.29'
$ws.Range('E35').Value = '  -5.08%  '
Set-TextValue 'D36' '1.00'
$ws.Range('E36').Value = '  +0.08%  '
Set-TextValue 'D37' '0.146'
$ws.Range('E37').Value = '  -3.67%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D38' '3.42'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D39' '35.03'
$ws.Range('E39').Value = '  -8.74%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D40' '0.0₃0733'
$ws.Range('E40').Value = '  -10.63%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D41' '3.148.38'
$ws.Range('E41').Value = '  +0.64%  '
Set-TextValue 'D42' '0.364'
$ws.Range('E42').Value = '  -6.29%  '
Set-TextValue 'D43' '0.997'
$ws.Range('E43').Value = '  -0.01%  '
Set-TextValue 'D44' '2.79'
$ws.Range('E44').Value = '  -2.70%  '
Set-TextValue 'D45' '3.17'
$ws.Range('E45').Value = '  -3.67%  '
Set-TextValue 'D46' '0.0399'
$ws.Range('E46').Value = '  -4.78%  '
Set-TextValue 'D47' '2.41'
$ws.Range('E47').Value = '  -8.37%  '
$ws.Range('E48').Value = '  -4.91%  '
Set-TextValue 'D49' '0.127'
$ws.Range('E49').Value = '  -5.30%  '
Set-TextValue 'D50' '132.58'
$ws.Range('E50').Value = '  -5.25%  '
Set-TextValue 'D51' '7.97'
$ws.Range('E51').Value = '  -6.94%  '
